# Ajout des fichiers générés utiles.
# Fill in the room/location ("salle") column (F) of the schedule with the
# values that were generated for each session.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F3").Value  = "U3-216"
$ws.Range("F6").Value  = "U3-308"
$ws.Range("F7").Value  = "U3-308"
$ws.Range("F10").Value = "U3-206"
$ws.Range("F13").Value = "U3-Amphi"
